$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1544.4576
$ws.Range("I40").Value = 1447.6086
$ws.Range("J40").Value = 1887.1538
$ws.Range("K40").Value = 1447.6086
$ws.Range("L40").Value = 1887.1538
$ws.Range("M40").Value = -1272.6086
$ws.Range("N40").Value = -2237.1538

$ws.Range("H62").Value = 14020.895
$ws.Range("I62").Value = 27323.875
$ws.Range("J62").Value = 4346
$ws.Range("K62").Value = 27323.875
$ws.Range("L62").Value = 4346
$ws.Range("M62").Value = -26699.875
$ws.Range("N62").Value = -5594

$ws.Range("H65").Value = 14020.895
$ws.Range("I65").Value = 27323.875
$ws.Range("J65").Value = 4346
$ws.Range("K65").Value = 136619.375
$ws.Range("L65").Value = 21730
$ws.Range("M65").Value = -133499.375
$ws.Range("N65").Value = -27970

$ws.Range("H94").Value = 1522.7778
$ws.Range("I94").Value = 1522.7778
$ws.Range("K94").Value = 1522.7778
$ws.Range("M94").Value = -1071.7778

$ws.Range("H113").Value = 3030.5898
$ws.Range("I113").Value = 2745.25
$ws.Range("J113").Value = 3157.4075
$ws.Range("K113").Value = 2745.25
$ws.Range("L113").Value = 3157.4075
$ws.Range("M113").Value = 508.75
$ws.Range("N113").Value = -9665.407499999999

$ws.Range("H116").Value = 59875.832
$ws.Range("I116").Value = 75868.92999999999
$ws.Range("K116").Value = 75868.92999999999
$ws.Range("M116").Value = -72426.92999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2776.5557
$ws.Range("I61").Value = 2044.4615
$ws.Range("J61").Value = 4680
$ws.Range("K61").Value = 2044.4615
$ws.Range("L61").Value = 4680
$ws.Range("M61").Value = -1832.4615
$ws.Range("N61").Value = -5104

$ws.Range("H101").Value = 31514
$ws.Range("J101").Value = 31514
$ws.Range("L101").Value = 31514
$ws.Range("N101").Value = -38004

$ws.Range("H110").Value = 1229.0555
$ws.Range("I110").Value = 1061.5333
$ws.Range("K110").Value = 1061.5333
$ws.Range("M110").Value = 983.4666999999999

$ws.Range("H132").Value = 3444.9768
$ws.Range("I132").Value = 3376.5
$ws.Range("K132").Value = 10129.5
$ws.Range("M132").Value = -7599.5

$ws.Range("H136").Value = 2776.5557
$ws.Range("I136").Value = 2044.4615
$ws.Range("J136").Value = 4680
$ws.Range("K136").Value = 6133.3845
$ws.Range("L136").Value = 14040
$ws.Range("M136").Value = -3583.3845
$ws.Range("N136").Value = -19140

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 26240.934
$ws.Range("I134").Value = 30259.5
$ws.Range("J134").Value = 10166.667
$ws.Range("K134").Value = 90778.5
$ws.Range("L134").Value = 30500.001
$ws.Range("M134").Value = -88243.5
$ws.Range("N134").Value = -35570.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2432.923
$ws.Range("I58").Value = 1200
$ws.Range("J58").Value = 2802.8
$ws.Range("K58").Value = 1200
$ws.Range("L58").Value = 2802.8
$ws.Range("M58").Value = -997
$ws.Range("N58").Value = -3208.8

$ws.Range("H99").Value = 55099.527
$ws.Range("I99").Value = 93007.09
$ws.Range("J99").Value = 2976.625
$ws.Range("K99").Value = 93007.09
$ws.Range("L99").Value = 2976.625
$ws.Range("M99").Value = -91509.09
$ws.Range("N99").Value = -5972.625

$ws.Range("H122").Value = 1725.96
$ws.Range("I122").Value = 2556.2727
$ws.Range("J122").Value = 1073.5714
$ws.Range("K122").Value = 7668.8181
$ws.Range("L122").Value = 3220.7142
$ws.Range("M122").Value = -5218.8181
$ws.Range("N122").Value = -8120.7142

$ws.Range("H126").Value = 55099.527
$ws.Range("I126").Value = 93007.09
$ws.Range("J126").Value = 2976.625
$ws.Range("K126").Value = 279021.27
$ws.Range("L126").Value = 8929.875
$ws.Range("M126").Value = -276551.27
$ws.Range("N126").Value = -13869.875

$ws.Range("H132").Value = 2060.1
$ws.Range("I132").Value = 1053.125
$ws.Range("J132").Value = 2989.6155
$ws.Range("K132").Value = 3159.375
$ws.Range("L132").Value = 8968.8465
$ws.Range("M132").Value = -629.375
$ws.Range("N132").Value = -14028.8465

$ws.Range("H134").Value = 1647.7675
$ws.Range("I134").Value = 945.9259
$ws.Range("J134").Value = 2832.125
$ws.Range("K134").Value = 2837.7777
$ws.Range("L134").Value = 8496.375
$ws.Range("M134").Value = -302.7776999999996
$ws.Range("N134").Value = -13566.375

$ws.Range("H136").Value = 2432.923
$ws.Range("I136").Value = 1200
$ws.Range("J136").Value = 2802.8
$ws.Range("K136").Value = 3600
$ws.Range("L136").Value = 8408.400000000001
$ws.Range("M136").Value = -1050
$ws.Range("N136").Value = -13508.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5355.5557
$ws.Range("I56").Value = 5355.5557
$ws.Range("K56").Value = 5355.5557
$ws.Range("M56").Value = -4825.5557

$ws.Range("H134").Value = 3121.2942
$ws.Range("I134").Value = 2883
$ws.Range("J134").Value = 3251.2727
$ws.Range("K134").Value = 8649
$ws.Range("L134").Value = 9753.8181
$ws.Range("M134").Value = -3579
$ws.Range("N134").Value = -19893.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3926.484
$ws.Range("I70").Value = 3590.0625
$ws.Range("J70").Value = 4285.3335
$ws.Range("K70").Value = 3590.0625
$ws.Range("L70").Value = 4285.3335
$ws.Range("M70").Value = -3320.0625
$ws.Range("N70").Value = -4825.3335

$ws.Range("H73").Value = 3926.484
$ws.Range("I73").Value = 3590.0625
$ws.Range("J73").Value = 4285.3335
$ws.Range("K73").Value = 3590.0625
$ws.Range("L73").Value = 4285.3335
$ws.Range("M73").Value = -2654.0625
$ws.Range("N73").Value = -6157.3335

$ws.Range("H132").Value = 3041.9592
$ws.Range("I132").Value = 2698.923
$ws.Range("J132").Value = 4379.8
$ws.Range("K132").Value = 8096.768999999999
$ws.Range("L132").Value = 13139.4
$ws.Range("M132").Value = -5566.768999999999
$ws.Range("N132").Value = -18199.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4460.5884
$ws.Range("I62").Value = 4990
$ws.Range("J62").Value = 4390
$ws.Range("K62").Value = 4990
$ws.Range("L62").Value = 4390
$ws.Range("M62").Value = -4366
$ws.Range("N62").Value = -5638

$ws.Range("H65").Value = 4460.5884
$ws.Range("I65").Value = 4990
$ws.Range("J65").Value = 4390
$ws.Range("K65").Value = 24950
$ws.Range("L65").Value = 21950
$ws.Range("M65").Value = -21830
$ws.Range("N65").Value = -28190

$ws.Range("H96").Value = 5620.4287
$ws.Range("I96").Value = 2380
$ws.Range("J96").Value = 13721.5
$ws.Range("K96").Value = 2380
$ws.Range("L96").Value = 13721.5
$ws.Range("M96").Value = -1007
$ws.Range("N96").Value = -16467.5
